# Weekly update: a new record is added at the top of the Papaya data block
# (row 54). All existing records from row 54 down shift one row down
# (54->55, ..., 87->88), and the former last record (old row 88) becomes
# the new last record at row 89. The new row 54 carries the new weekly
# price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 54; this pushes rows 54-88 down to 55-89
# and keeps all their existing values/styles intact.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with this week's new observation.
$ws.Range("A54").Value = 10
$ws.Range("B54").Value = 'Vega Modelo de Temuco'
$ws.Range("C54").Value = 'La Araucanía'
$ws.Range("D54").Value = 44977
$ws.Range("E54").Value = 9
$ws.Range("F54").Value = 'Fruta'
$ws.Range("G54").Value = 100108
$ws.Range("H54").Value = 'Tropicales y subtropicales'
$ws.Range("I54").Value = 100108004
$ws.Range("J54").Value = 'Papaya'
$ws.Range("K54").Value = 'Cultivar IV Región'
$ws.Range("L54").Value = 'Primera'
$ws.Range("M54").Value = 50
$ws.Range("N54").Value = 40000
$ws.Range("O54").Value = 40000
$ws.Range("P54").Value = 40000
$ws.Range("Q54").Value = '$/caja 15 kilos granel'
$ws.Range("R54").Value = 'Provincia del Elquí'
$ws.Range("S54").Value = 2667
$ws.Range("T54").Value = 15
